# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (H) season totals updated
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 393
$wsOff.Range("C2").Value = 256
$wsOff.Range("D2").Value = 56
$wsOff.Range("E2").Value = 24
$wsOff.Range("F2").Value = 9
$wsOff.Range("G2").Value = 4

# DEF sheet - row 2 (H) season totals updated
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 430
$wsDef.Range("C2").Value = 305
$wsDef.Range("D2").Value = 87
